$wb = $excel.ActiveWorkbook

# --- Sheet "DQ_Report" ---
$wsReport = $wb.Worksheets.Item("DQ_Report")

$wsReport.Range("D2").Value = "Kodierung ist nicht eindeutig. Relation G70 - 586 ist im BfArM nicht vorhanden. "
$wsReport.Range("D3").Value = "Kodierung ist nicht eindeutig. Relation G70 - 589 ist im BfArM nicht vorhanden. "
$wsReport.Range("D13").Value = "ICD10 Kodierung E66.89 ist nicht eindeutig. ICD10-Orpha Relation ist gemäß Tracer-Diagnosenliste vom Typ 1-m.  Fehlendes Orpha_Kode.  "
$wsReport.Range("D14").Value = "ICD10 Kodierung E75.2 ist nicht eindeutig. ICD10-Orpha Relation ist gemäß Tracer-Diagnosenliste vom Typ 1-m.  Fehlendes Orpha_Kode.  "

# --- Sheet "Statistik" ---
$wsStat = $wb.Worksheets.Item("Statistik")
$wsStat.Range("E2").Value = 97.09999999999999
